$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 585
$ws.Range("F6").Value = 28
$ws.Range("F8").Value = 5308
$ws.Range("F9").Value = 1488
$ws.Range("F11").Value = 3079
$ws.Range("F14").Value = 1279
$ws.Range("F15").Value = 4221
$ws.Range("F17").Value = 885
$ws.Range("F19").Value = 2596
$ws.Range("F20").Value = 29
$ws.Range("F21").Value = 21
$ws.Range("F24").Value = 965
$ws.Range("F29").Value = 1076
$ws.Range("F30").Value = 366
$ws.Range("F31").Value = 36
$ws.Range("F32").Value = 125
$ws.Range("F34").Value = 250
$ws.Range("F36").Value = 2163
$ws.Range("F38").Value = 35
$ws.Range("F41").Value = 271
$ws.Range("F43").Value = 649
$ws.Range("F44").Value = 395
$ws.Range("F45").Value = 318
$ws.Range("F46").Value = 206

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 731

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 731
$ws.Range("F4").Value = 585
$ws.Range("F5").Value = 28
$ws.Range("F7").Value = 5308
$ws.Range("F8").Value = 1488
$ws.Range("F11").Value = 3079
$ws.Range("F13").Value = 1279
$ws.Range("F14").Value = 4221
$ws.Range("F18").Value = 2596
$ws.Range("F20").Value = 29
$ws.Range("F23").Value = 21
$ws.Range("F26").Value = 966
$ws.Range("F32").Value = 1076
$ws.Range("F33").Value = 366
$ws.Range("F34").Value = 36
$ws.Range("F36").Value = 2163
$ws.Range("F38").Value = 35
$ws.Range("F43").Value = 271
$ws.Range("F44").Value = 649
$ws.Range("F45").Value = 395
$ws.Range("F46").Value = 318
$ws.Range("F47").Value = 206
